# Convert the numeric codes in column H ("TipoFort") into their text
# labels (PHM / PH / P3 / P), matching the new lookup strings added to
# the shared-string table, and update the sheet's saved selection.
#
# The shared-string table must end up with the new strings appended in
# the order PHM, PH, P3, P (right after the existing 8 strings), so we
# assign the values to the affected cells grouped by label, in that
# exact order, rather than sweeping top-to-bottom through the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# old value 1 -> "PHM"
$r = $ws.Range("H16,H22,H24,H66,H68,H73,H83,H86:H87")
foreach ($area in $r.Areas) { $area.Value = "PHM" }

# old value 2 -> "PH" (except H10:H11, see below)
$r = $ws.Range("H8:H9,H12:H15,H17:H21,H23,H25:H39,H52:H54,H56:H65,H67,H69:H70,H72,H74:H82,H84:H85,H88:H91")
foreach ($area in $r.Areas) { $area.Value = "PH" }

# old value 3 -> "P3" (plus H10:H11, which were old value 2)
$r = $ws.Range("H2:H5,H10:H11,H40:H51,H55,H71,H92")
foreach ($area in $r.Areas) { $area.Value = "P3" }

# old value 4 -> "P"
$r = $ws.Range("H6:H7")
foreach ($area in $r.Areas) { $area.Value = "P" }

# Update the sheet view: clear the scrolled-down top-left cell and move
# the active selection to I12.
[void]$ws.Range("I12").Select()
